$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

# --- Column D (Price) text updates ---
Set-TextValue "D2" "29.390.67"
Set-TextValue "D3" "1.846.72"
Set-TextValue "D4" "0.9987"
Set-TextValue "D5" "240.16"
Set-TextValue "D6" "0.6314"
Set-TextValue "D7" "1.000"
Set-TextValue "D8" "0.07548"
Set-TextValue "D9" "0.2964"
Set-TextValue "D10" "24.59"
Set-TextValue "D11" "0.07728"
Set-TextValue "D12" "1.846.04"
Set-TextValue "D13" "4.999"
Set-TextValue "D14" "0.6849"
Set-TextValue "D15" "0.00001002"
Set-TextValue "D16" "83.10"
Set-TextValue "D17" "6.181"
Set-TextValue "D18" "29.415.97"
Set-TextValue "D19" "229.96"
Set-TextValue "D21" "0.9995"
Set-TextValue "D22" "7.573"
Set-TextValue "D23" "1.001"
Set-TextValue "D24" "157.05"
Set-TextValue "D25" "0.1399"
Set-TextValue "D28" "1.464"
Set-TextValue "D29" "0.05737"
Set-TextValue "D30" "1.250"
Set-TextValue "D32" "4.037"
Set-TextValue "D33" "1.855"
Set-TextValue "D35" "0.7173"
Set-TextValue "D36" "2.592"
Set-TextValue "D37" "1.251.37"
Set-TextValue "D38" "0.01818"
Set-TextValue "D39" "2.787"
Set-TextValue "D40" "6.212"
Set-TextValue "D41" "0.9078"
Set-TextValue "D43" "1.993.18"
Set-TextValue "D44" "101.78"
Set-TextValue "D45" "66.42"

# --- Column E (Volume/1h) text updates ---
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -0.29%  "
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("E16").Value = "  -0.75%  "
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("E19").Value = "  -1.77%  "
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -0.39%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("E24").Value = "  +0.81%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  -0.48%  "
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("E28").Value = "  -1.13%  "
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("E30").Value = "  -2.71%  "
$ws.Range("E32").Value = "  -0.08%  "
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("E34").Value = "  -1.13%  "
$ws.Range("E37").Value = "  +0.95%  "
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  -0.38%  "
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("E43").Value = "  -2.43%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  -1.53%  "

# --- Rows 46-51: new coin inserted (BabyDogeCoin), shifting rows down ---
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D46" "0.00000000118"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.064"
$ws.Range("E47").Value = "  -2.87%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.170"
$ws.Range("E48").Value = "  +0.08%  "
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D49" "0.4030"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D50" "1.709"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D51" "0.1131"
$ws.Range("E51").Value = "  +0.90%  "
